$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F10"
$ws.Range("C2").Value = "F3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 1.288450333333333
$ws.Range("H2").Value = 3.865351
$ws.Range("I2").Value = 0.983869567795948
$ws.Range("J2").Value = 0.9838695677959479
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 198.2465873333333
$ws.Range("N2").Value = 594.7397619999999
$ws.Range("O2").Value = 0.9851515664921635
$ws.Range("P2").Value = 0.9851515664921635
$ws.Range("Q2").Value = 255.4308815318291
$ws.Range("R2").Value = 2298.877933786462
$ws.Range("S2").Value = 0.9692606459381461
$ws.Range("T2").Value = 0.969260645938146

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F10"
$ws.Range("C3").Value = "F3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 1.288450333333333
$ws.Range("H3").Value = 3.865351
$ws.Range("I3").Value = 0.983869567795948
$ws.Range("J3").Value = 0.9838695677959479
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 2.988018666666667
$ws.Range("N3").Value = 8.964056
$ws.Range("O3").Value = 0.01484843350783645
$ws.Range("P3").Value = 0.01484843350783645
$ws.Range("Q3").Value = 3.849913647072889
$ws.Range("R3").Value = 34.649222823656
$ws.Range("S3").Value = 0.01460892185780193
$ws.Range("T3").Value = 0.01460892185780192

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F10"
$ws.Range("C4").Value = "F3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01545266666666667
$ws.Range("H4").Value = 0.046358
$ws.Range("I4").Value = 0.0117997629255104
$ws.Range("J4").Value = 0.0117997629255104
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 198.2465873333333
$ws.Range("N4").Value = 594.7397619999999
$ws.Range("O4").Value = 0.9851515664921635
$ws.Range("P4").Value = 0.9851515664921635
$ws.Range("Q4").Value = 3.063438431866222
$ws.Range("R4").Value = 27.570945886796
$ws.Range("S4").Value = 0.01162455493030273
$ws.Range("T4").Value = 0.01162455493030272

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F10"
$ws.Range("C5").Value = "F3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1.0
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01545266666666667
$ws.Range("H5").Value = 0.046358
$ws.Range("I5").Value = 0.0117997629255104
$ws.Range("J5").Value = 0.0117997629255104
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 2.988018666666667
$ws.Range("N5").Value = 8.964056
$ws.Range("O5").Value = 0.01484843350783645
$ws.Range("P5").Value = 0.01484843350783645
$ws.Range("Q5").Value = 0.04617285644977778
$ws.Range("R5").Value = 0.415555708048
$ws.Range("S5").Value = 0.000175207995207675
$ws.Range("T5").Value = 0.0001752079952076749

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F10"
$ws.Range("C6").Value = "F3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1.0
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.005671333333333334
$ws.Range("H6").Value = 0.017014
$ws.Range("I6").Value = 0.004330669278541654
$ws.Range("J6").Value = 0.004330669278541653
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 198.2465873333333
$ws.Range("N6").Value = 594.7397619999999
$ws.Range("O6").Value = 0.9851515664921635
$ws.Range("P6").Value = 0.9851515664921635
$ws.Range("Q6").Value = 1.124322478963111
$ws.Range("R6").Value = 10.118902310668
$ws.Range("S6").Value = 0.004266365623714798
$ws.Range("T6").Value = 0.004266365623714797

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F10"
$ws.Range("C7").Value = "F3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1.0
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.005671333333333334
$ws.Range("H7").Value = 0.017014
$ws.Range("I7").Value = 0.004330669278541654
$ws.Range("J7").Value = 0.004330669278541653
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 2.988018666666667
$ws.Range("N7").Value = 8.964056
$ws.Range("O7").Value = 0.01484843350783645
$ws.Range("P7").Value = 0.01484843350783645
$ws.Range("Q7").Value = 0.01694604986488889
$ws.Range("R7").Value = 0.152514448784
$ws.Range("S7").Value = 0.00006430365482685582
$ws.Range("T7").Value = 0.0000643036548268558

$ws.Range("A8:T10").Delete()
